# Update the "Förändrad" (Changed) date column (C) for rows 2-118
# from serial date 45179 (2023-09-10) to 45180 (2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($row in 2..118) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
